$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<he>"
$ws.Range("C2").Value = 58

$ws.Range("B3").Value = "<may>"
$ws.Range("C3").Value = 58

$ws.Range("B4").Value = "<who>"

$ws.Range("B5").Value = "<water>"
$ws.Range("C5").Value = 60

$ws.Range("B6").Value = "<so>"

$ws.Range("B7").Value = "<these>"
$ws.Range("C7").Value = 59

$ws.Range("B8").Value = "<her>"
$ws.Range("C8").Value = 59

$ws.Range("B9").Value = "<them>"
$ws.Range("C9").Value = 59

$ws.Range("B10").Value = "<water>"

$ws.Range("B11").Value = "<which>"
$ws.Range("C11").Value = 59

$ws.Range("B12").Value = "<is>"
$ws.Range("C12").Value = 60

$ws.Range("B13").Value = "<his>"

$ws.Range("B14").Value = "<a>"
